$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.159.18'
$ws.Range("E2").Value = '  +0.44%  '
$ws.Range("D3").Value = '1.830.74'
$ws.Range("E3").Value = '  -0.13%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9992'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '241.56'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.84%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6584'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.84%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9999'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07400'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.17%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2923'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.87%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '22.87'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.82%  '
$ws.Range("E11").Value = '  +1.37%  '
$ws.Range("D12").Value = '1.837.82'
$ws.Range("E12").Value = '  +0.37%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.992'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.18%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6650'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.06%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '82.76'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -4.06%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.107'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.41%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000008431'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.59%  '
$ws.Range("D18").Value = '29.160.89'
$ws.Range("E18").Value = '  +0.48%  '
$ws.Range("D19").Value = '2.081.69'
$ws.Range("E19").Value = '  +0.37%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '226.89'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.10%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.43'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.10%  '
$ws.Range("E22").Value = '  +0.13%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.119'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.60%  '
$ws.Range("E24").Value = '  +0.05%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '158.77'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.81%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.597'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.76%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1390'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.64%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '17.91'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.17%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.517'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.12%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.112'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.82%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.041'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.69%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.187'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.72%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05249'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.95%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.863'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.60%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7395'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.29%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.141'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.64%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.655'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.08%  '
$ws.Range("D38").Value = '1.297.94'
$ws.Range("E38").Value = '  +0.42%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01788'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.90%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.732'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.92%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9192'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.35%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.953'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.56%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.08500'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.89%  '
$ws.Range("E44").Value = '  +0.03%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '102.24'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.81%  '
$ws.Range("D46").Value = '1.976.30'
$ws.Range("E46").Value = '  +0.18%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5142'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.63%  '
$ws.Range("B48").Value = 'RenderToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.750'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.07%  '
$ws.Range("B49").Value = 'BabyDogeCoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.00000000120'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -8.29%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '63.31'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.33%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05841'
$ws.Range("D51").Style = "Normal"
